# Updated cryptos list with GitHub Actions: refresh Price (column D) and
# Volume(1h) (column E) figures, and swap the Arweave / InjectiveProtocol
# rows (50/51) to reflect their new ranking order.
#
# Note: several Price values are plain decimals (e.g. "564.53") that Excel
# would otherwise auto-convert to a number. We write them with a leading
# apostrophe (forces text entry) and then reset the cell style back to
# "Normal" so the stored value is text without leaving a stray
# "quote prefix" number format applied to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.298.08"
$ws.Range("E2").Value = "  +4.20%  "
$ws.Range("D3").Value = "2.995.19"
$ws.Range("E3").Value = "  +3.88%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'564.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.41%  "
$ws.Range("D6").Value = "'138.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +12.55%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.521"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.38%  "
$ws.Range("D9").Value = "2.988.40"
$ws.Range("E9").Value = "  +3.76%  "
$ws.Range("E10").Value = "  +10.33%  "
$ws.Range("D11").Value = "'5.03"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.97%  "
$ws.Range("E12").Value = "  +5.03%  "
$ws.Range("E13").Value = "  +10.27%  "
$ws.Range("E14").Value = "  +4.25%  "
$ws.Range("E15").Value = "  +2.94%  "
$ws.Range("D16").Value = "3.491.40"
$ws.Range("E16").Value = "  +3.93%  "
$ws.Range("E17").Value = "  +7.69%  "
$ws.Range("D18").Value = "2.995.83"
$ws.Range("E18").Value = "  +3.94%  "
$ws.Range("D19").Value = "59.230.80"
$ws.Range("E19").Value = "  +3.92%  "
$ws.Range("D20").Value = "'430.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.68%  "
$ws.Range("D21").Value = "'13.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.66%  "
$ws.Range("D22").Value = "'0.714"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.57%  "
$ws.Range("E23").Value = "  +5.63%  "
$ws.Range("E24").Value = "  +6.85%  "
$ws.Range("D25").Value = "'80.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("E28").Value = "  +11.42%  "
$ws.Range("D29").Value = "'2.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.25%  "
$ws.Range("E30").Value = "  +8.87%  "
$ws.Range("E31").Value = "  +4.60%  "
$ws.Range("D32").Value = "'6.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.14%  "
$ws.Range("D33").Value = "'0.0992"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.93%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.59%  "
$ws.Range("E35").Value = "  +24.02%  "
$ws.Range("D36").Value = "'5.79"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.64%  "
$ws.Range("E37").Value = "  +4.09%  "
$ws.Range("E38").Value = "  +3.06%  "
$ws.Range("D39").Value = "'8.66"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.72%  "
$ws.Range("E40").Value = "  +15.72%  "
$ws.Range("D41").Value = "'402.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +12.62%  "
$ws.Range("D42").Value = "'0.0352"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.68%  "
$ws.Range("D43").Value = "2.756.25"
$ws.Range("E43").Value = "  +4.95%  "
$ws.Range("E44").Value = "  +3.27%  "
$ws.Range("E45").Value = "  +10.59%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").Value = "'123.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.56%  "
$ws.Range("D48").Value = "'2.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.68%  "
$ws.Range("E49").Value = "  +2.72%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'23.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.87%  "
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").Value = "'32.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +19.82%  "
